$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.031.80"
$ws.Range("E2").Value = "  +0.61%  "

# Row 3
$ws.Range("D3").Value = "1.643.73"
$ws.Range("E3").Value = "  +0.67%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.37"
$ws.Range("E5").Value = "  +0.84%  "

# Row 7
$ws.Range("E7").Value = "  +0.14%  "

# Row 8
$ws.Range("E8").Value = "  +0.50%  "

# Row 9
$ws.Range("E9").Value = "  +1.28%  "

# Row 10
$ws.Range("E10").Value = "  +0.65%  "

# Row 11
$ws.Range("E11").Value = "  +0.63%  "

# Row 12
$ws.Range("D12").Value = "1.871.83"
$ws.Range("E12").Value = "  +0.65%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.674.97"
$ws.Range("E13").Value = "  +2.74%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.28"
$ws.Range("E14").Value = "  +0.86%  "

# Row 15
$ws.Range("E15").Value = "  +0.24%  "

# Row 16
$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("E16").Value = "  +1.39%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.41"
$ws.Range("E17").Value = "  +1.02%  "

# Row 18
$ws.Range("D18").Value = "26.050.35"
$ws.Range("E18").Value = "  +0.55%  "

# Row 19
$ws.Range("E19").Value = "  +0.28%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "195.75"
$ws.Range("E20").Value = "  +1.43%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.36"
$ws.Range("E21").Value = "  -0.29%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.94"
$ws.Range("E22").Value = "  +0.04%  "

# Row 24
$ws.Range("E24").Value = "  +5.11%  "

# Row 25
$ws.Range("E25").Value = "  +0.23%  "

# Row 26
$ws.Range("E26").Value = "  +0.72%  "

# Row 27
$ws.Range("E27").Value = "  +0.28%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.91"
$ws.Range("E28").Value = "  +0.95%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  +0.71%  "

# Row 30
$ws.Range("E30").Value = "  +1.31%  "

# Row 31
$ws.Range("E31").Value = "  +0.20%  "

# Row 32
$ws.Range("E32").Value = "  -0.01%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.26"
$ws.Range("E33").Value = "  +1.41%  "

# Row 34
$ws.Range("E34").Value = "  -2.86%  "

# Row 35
$ws.Range("E35").Value = "  +1.18%  "

# Row 36
$ws.Range("E36").Value = "  +0.62%  "

# Row 37
$ws.Range("D37").Value = "1.135.12"
$ws.Range("E37").Value = "  -0.21%  "

# Row 38
$ws.Range("E38").Value = "  -1.37%  "

# Row 39
$ws.Range("E39").Value = "  -0.67%  "

# Row 40
$ws.Range("E40").Value = "  +0.32%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.49"
$ws.Range("E41").Value = "  +0.85%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.43"
$ws.Range("E42").Value = "  +0.27%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.798"
$ws.Range("E43").Value = "  -0.70%  "

# Row 44
$ws.Range("D44").Value = "1.781.33"

# Row 45
$ws.Range("E45").Value = "  +5.85%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.83"
$ws.Range("E46").Value = "  +1.10%  "

# Row 47
$ws.Range("E47").Value = "  -0.81%  "

# Row 48
$ws.Range("E48").Value = "  +0.84%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.79"
$ws.Range("E49").Value = "  +2.11%  "

# Row 50
$ws.Range("E50").Value = "  -0.35%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0961"
$ws.Range("E51").Value = "  +0.02%  "
